$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 635 (the "「私はサッカーが好きです」" post),
# which shifts all subsequent rows up by one.
$ws.Rows.Item(635).Delete()
